$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.176022
$ws.Range("H2").Value = 0.5280659999999999
$ws.Range("I2").Value = 0.03293066697281707
$ws.Range("J2").Value = 0.03293066697281707
$ws.Range("M2").Value = 115.2213693333333
$ws.Range("N2").Value = 345.664108
$ws.Range("O2").Value = 0.2787408744545015
$ws.Range("P2").Value = 0.2787408744545015
$ws.Range("Q2").Value = 20.281495872792
$ws.Range("R2").Value = 182.533462855128
$ws.Range("S2").Value = 0.009179122908373003
$ws.Range("T2").Value = 0.009179122908373

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.176022
$ws.Range("H3").Value = 0.5280659999999999
$ws.Range("I3").Value = 0.03293066697281707
$ws.Range("J3").Value = 0.03293066697281707
$ws.Range("O3").Value = 0.44716501655323
$ws.Range("P3").Value = 0.4471650165532299
$ws.Range("Q3").Value = 32.536223671644
$ws.Range("R3").Value = 292.826013044796
$ws.Range("S3").Value = 0.01472544224200865
$ws.Range("T3").Value = 0.01472544224200865

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.176022
$ws.Range("H4").Value = 0.5280659999999999
$ws.Range("I4").Value = 0.03293066697281707
$ws.Range("J4").Value = 0.03293066697281707
$ws.Range("M4").Value = 60.55095666666667
$ws.Range("N4").Value = 181.65287
$ws.Range("O4").Value = 0.1464834753134679
$ws.Range("P4").Value = 0.1464834753134678
$ws.Range("Q4").Value = 10.65830049438
$ws.Range("R4").Value = 95.92470444941999
$ws.Range("S4").Value = 0.004823798542568682
$ws.Range("T4").Value = 0.00482379854256868

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.176022
$ws.Range("H5").Value = 0.5280659999999999
$ws.Range("I5").Value = 0.03293066697281707
$ws.Range("J5").Value = 0.03293066697281707
$ws.Range("M5").Value = 52.74960833333333
$ws.Range("N5").Value = 158.248825
$ws.Range("O5").Value = 0.1276106336788006
$ws.Range("P5").Value = 0.1276106336788006
$ws.Range("Q5").Value = 9.285091558049999
$ws.Range("R5").Value = 83.56582402244999
$ws.Range("S5").Value = 0.004202303279866738
$ws.Range("T5").Value = 0.004202303279866738

$ws.Range("I6").Value = 0.8002039325901205
$ws.Range("J6").Value = 0.8002039325901203
$ws.Range("M6").Value = 115.2213693333333
$ws.Range("N6").Value = 345.664108
$ws.Range("O6").Value = 0.2787408744545015
$ws.Range("P6").Value = 0.2787408744545015
$ws.Range("Q6").Value = 492.8334057009871
$ws.Range("R6").Value = 4435.500651308884
$ws.Range("S6").Value = 0.2230495439121012
$ws.Range("T6").Value = 0.2230495439121011

$ws.Range("I7").Value = 0.8002039325901205
$ws.Range("J7").Value = 0.8002039325901203
$ws.Range("O7").Value = 0.44716501655323
$ws.Range("P7").Value = 0.4471650165532299
$ws.Range("S7").Value = 0.357823204762621
$ws.Range("T7").Value = 0.3578232047626209

$ws.Range("I8").Value = 0.8002039325901205
$ws.Range("J8").Value = 0.8002039325901203
$ws.Range("M8").Value = 60.55095666666667
$ws.Range("N8").Value = 181.65287
$ws.Range("O8").Value = 0.1464834753134679
$ws.Range("P8").Value = 0.1464834753134678
$ws.Range("Q8").Value = 258.9930528091123
$ws.Range("R8").Value = 2330.93747528201
$ws.Range("S8").Value = 0.1172166530053048
$ws.Range("T8").Value = 0.1172166530053048

$ws.Range("I9").Value = 0.8002039325901205
$ws.Range("J9").Value = 0.8002039325901203
$ws.Range("M9").Value = 52.74960833333333
$ws.Range("N9").Value = 158.248825
$ws.Range("O9").Value = 0.1276106336788006
$ws.Range("P9").Value = 0.1276106336788006
$ws.Range("Q9").Value = 225.6245458175528
$ws.Range("R9").Value = 2030.620912357975
$ws.Range("S9").Value = 0.1021145309100935
$ws.Range("T9").Value = 0.1021145309100935

$ws.Range("G10").Value = 0.891934
$ws.Range("H10").Value = 2.675802
$ws.Range("I10").Value = 0.1668654004370625
$ws.Range("J10").Value = 0.1668654004370625
$ws.Range("M10").Value = 115.2213693333333
$ws.Range("N10").Value = 345.664108
$ws.Range("O10").Value = 0.2787408744545015
$ws.Range("P10").Value = 0.2787408744545015
$ws.Range("Q10").Value = 102.7698568349573
$ws.Range("R10").Value = 924.928711514616
$ws.Range("S10").Value = 0.04651220763402737
$ws.Range("T10").Value = 0.04651220763402737

$ws.Range("G11").Value = 0.891934
$ws.Range("H11").Value = 2.675802
$ws.Range("I11").Value = 0.1668654004370625
$ws.Range("J11").Value = 0.1668654004370625
$ws.Range("O11").Value = 0.44716501655323
$ws.Range("P11").Value = 0.4471650165532299
$ws.Range("Q11").Value = 164.866687825068
$ws.Range("R11").Value = 1483.800190425612
$ws.Range("S11").Value = 0.07461636954860043
$ws.Range("T11").Value = 0.07461636954860042

$ws.Range("G12").Value = 0.891934
$ws.Range("H12").Value = 2.675802
$ws.Range("I12").Value = 0.1668654004370625
$ws.Range("J12").Value = 0.1668654004370625
$ws.Range("M12").Value = 60.55095666666667
$ws.Range("N12").Value = 181.65287
$ws.Range("O12").Value = 0.1464834753134679
$ws.Range("P12").Value = 0.1464834753134678
$ws.Range("Q12").Value = 54.00745698352667
$ws.Range("R12").Value = 486.06711285174
$ws.Range("S12").Value = 0.02444302376559438
$ws.Range("T12").Value = 0.02444302376559438

$ws.Range("G13").Value = 0.891934
$ws.Range("H13").Value = 2.675802
$ws.Range("I13").Value = 0.1668654004370625
$ws.Range("J13").Value = 0.1668654004370625
$ws.Range("M13").Value = 52.74960833333333
$ws.Range("N13").Value = 158.248825
$ws.Range("O13").Value = 0.1276106336788006
$ws.Range("P13").Value = 0.1276106336788006
$ws.Range("Q13").Value = 47.04916915918334
$ws.Range("R13").Value = 423.44252243265
$ws.Range("S13").Value = 0.02129379948884037
$ws.Range("T13").Value = 0.02129379948884037
